$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 396. This pushes the existing rows 396..421
# down to 397..422, preserving all of their data intact.
$ws.Rows.Item(396).Insert()

# Populate the newly inserted row 396 with its data. Columns A-C, E-J, N-R
# mirror the row that used to be at 396 (now at 397); D, K, L, M and P carry
# the new values for this entry.
$ws.Cells.Item(396, 1).Value = 5
$ws.Cells.Item(396, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(396, 3).Value = "Maule"
$ws.Cells.Item(396, 4).Value = 44706
$ws.Cells.Item(396, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(396, 5).Value = 7
$ws.Cells.Item(396, 6).Value = 100112043
$ws.Cells.Item(396, 7).Value = "Pepino ensalada"
$ws.Cells.Item(396, 8).Value = "Sin especificar"
$ws.Cells.Item(396, 9).Value = "Primera"
$ws.Cells.Item(396, 10).Value = 300
$ws.Cells.Item(396, 11).Value = 20000
$ws.Cells.Item(396, 12).Value = 20000
$ws.Cells.Item(396, 13).Value = 20000
$ws.Cells.Item(396, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(396, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(396, 16).Value = 333
$ws.Cells.Item(396, 17).Value = 60
$ws.Cells.Item(396, 18).Value = "Hortaliza"
